# Upload Y4_B2526_General_&_Special_surgery_1_reference_data.xlsx via attendance app
# Appends 4 new student rows (223003-223006) to the reference data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: Student ID, Name, Year, Group, Source File
$newRows = @(
    @{ Row = 312; Id = 223003; Name = "تسنيم محمد انور" },
    @{ Row = 313; Id = 223004; Name = "هاجر عماد حسين حسين" },
    @{ Row = 314; Id = 223005; Name = "عمر محمد احمد على محمد حفناوى" },
    @{ Row = 315; Id = 223006; Name = "محمد فتحى احمد الحسينى" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.Id
    $ws.Range("B$row").Value = $r.Name
    $ws.Range("C$row").Value = "Year 4"
    $ws.Range("D$row").Value = "B2F"
    $ws.Range("E$row").Value = "الغياب.xlsx"
}

# Match the existing formatting: every new row uses the "plain" style of row
# 311, except column B on rows 312 and 314 which keeps the "shaded" style
# used on row 310 (matching the source workbook's (slightly irregular)
# banding for the pasted-in rows).
$ws.Range("A311:E311").Copy() | Out-Null
$ws.Range("A312:E312").PasteSpecial(-4122) | Out-Null
$ws.Range("A313:E313").PasteSpecial(-4122) | Out-Null
$ws.Range("A314:E314").PasteSpecial(-4122) | Out-Null
$ws.Range("A315:E315").PasteSpecial(-4122) | Out-Null

$ws.Range("B310").Copy() | Out-Null
$ws.Range("B312").PasteSpecial(-4122) | Out-Null
$ws.Range("B314").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Reflect the new used range and move the on-screen selection to where the
# user ended up after pasting in the new rows.
$ws.Range("C316").Select() | Out-Null
